# Apply "Add data for 2022-12-15" update:
#  - Bump the "through December 06" labels to "through December 07"
#    in both the sheet name and the B1 header text.
#  - Update/insert the carjacking counts for several neighborhoods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab (Through 2022-12-06 -> Through 2022-12-07)
$ws.Name = "Through 2022-12-07"

# Update header text in B1
$ws.Range("B1").Value = "December 2022 (through December 07)"

# Updated / new cell values
$ws.Range("Z2").Value   = 8
$ws.Range("AX3").Value  = 2
$ws.Range("N4").Value   = 3
$ws.Range("CH6").Value  = 1
$ws.Range("AL7").Value  = 1
$ws.Range("CH10").Value = 1
$ws.Range("AX14").Value = 2
$ws.Range("BV14").Value = 3
$ws.Range("CH14").Value = 1
$ws.Range("N19").Value  = 1
$ws.Range("N20").Value  = 2
$ws.Range("Z20").Value  = 5
$ws.Range("BV20").Value = 2
$ws.Range("N22").Value  = 2
$ws.Range("BJ23").Value = 1
$ws.Range("BJ28").Value = 3
$ws.Range("AX37").Value = 1
$ws.Range("AL44").Value = 1
$ws.Range("N45").Value  = 2
$ws.Range("N47").Value  = 1
$ws.Range("Z49").Value  = 1
$ws.Range("N57").Value  = 2
$ws.Range("B64").Value  = 3
$ws.Range("BV65").Value = 1
$ws.Range("N99").Value  = 1
